# Auto-generated edit script: update Leve profit calculation cells across all 8 sheets
# per the Gilgamesh_Profits workbook diff (price refresh from scheduled runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 224.5  # H6
$ws.Cells.Item(6, 9).Value = 224.5  # I6
$ws.Cells.Item(6, 11).Value = 673.5  # K6
$ws.Cells.Item(6, 13).Value = -561.5  # M6

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 3.2  # H8
$ws.Cells.Item(8, 9).Value = 3.2222223  # I8
$ws.Cells.Item(8, 10).Value = 3.0  # J8
$ws.Cells.Item(8, 11).Value = 9.6666669  # K8
$ws.Cells.Item(8, 12).Value = 9.0  # L8
$ws.Cells.Item(8, 13).Value = 129.3333331  # M8
$ws.Cells.Item(8, 14).Value = -287.0  # N8

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 1909.0  # H43
$ws.Cells.Item(43, 10).Value = 1928.4286  # J43
$ws.Cells.Item(43, 12).Value = 1928.4286  # L43
$ws.Cells.Item(43, 14).Value = -2066.4286  # N43

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(96, 8).Value = 600.25  # H96
$ws.Cells.Item(96, 9).Value = 257.5  # I96
$ws.Cells.Item(96, 11).Value = 772.5  # K96
$ws.Cells.Item(96, 13).Value = 600.5  # M96

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 2817.6191  # H137
$ws.Cells.Item(137, 9).Value = 1848.75  # I137
$ws.Cells.Item(137, 11).Value = 5546.25  # K137
$ws.Cells.Item(137, 13).Value = -2996.25  # M137

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 5066.3335  # H138
$ws.Cells.Item(138, 9).Value = 4000.0  # I138
$ws.Cells.Item(138, 10).Value = 7199.0  # J138
$ws.Cells.Item(138, 11).Value = 12000.0  # K138
$ws.Cells.Item(138, 12).Value = 21597.0  # L138
$ws.Cells.Item(138, 13).Value = -6860.0  # M138
$ws.Cells.Item(138, 14).Value = -31877.0  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2180.1667  # H2
$ws.Cells.Item(2, 9).Value = 2016.3  # I2
$ws.Cells.Item(2, 10).Value = 2999.5  # J2
$ws.Cells.Item(2, 11).Value = 2016.3  # K2
$ws.Cells.Item(2, 12).Value = 2999.5  # L2
$ws.Cells.Item(2, 13).Value = -1903.3  # M2
$ws.Cells.Item(2, 14).Value = -3225.5  # N2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1389276.4  # H32
$ws.Cells.Item(32, 9).Value = 651295.25  # I32
$ws.Cells.Item(32, 11).Value = 651295.25  # K32
$ws.Cells.Item(32, 13).Value = -651008.25  # M32

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 893.8182  # H97
$ws.Cells.Item(97, 9).Value = 763.5  # I97
$ws.Cells.Item(97, 11).Value = 763.5  # K97
$ws.Cells.Item(97, 13).Value = -267.5  # M97

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 2180.1667  # H116
$ws.Cells.Item(116, 9).Value = 2016.3  # I116
$ws.Cells.Item(116, 10).Value = 2999.5  # J116
$ws.Cells.Item(116, 11).Value = 2016.3  # K116
$ws.Cells.Item(116, 12).Value = 2999.5  # L116
$ws.Cells.Item(116, 13).Value = 277.7  # M116
$ws.Cells.Item(116, 14).Value = -7587.5  # N116

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1398.0952  # H132
$ws.Cells.Item(132, 9).Value = 1092.1034  # I132
$ws.Cells.Item(132, 11).Value = 3276.3102  # K132
$ws.Cells.Item(132, 13).Value = -746.3101999999999  # M132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2180.1667  # H3
$ws.Cells.Item(3, 9).Value = 2016.3  # I3
$ws.Cells.Item(3, 10).Value = 2999.5  # J3
$ws.Cells.Item(3, 11).Value = 2016.3  # K3
$ws.Cells.Item(3, 12).Value = 2999.5  # L3
$ws.Cells.Item(3, 13).Value = -1902.3  # M3
$ws.Cells.Item(3, 14).Value = -3227.5  # N3

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4715.4614  # H86
$ws.Cells.Item(86, 9).Value = 4536.4546  # I86
$ws.Cells.Item(86, 11).Value = 4536.4546  # K86
$ws.Cells.Item(86, 13).Value = -3413.4546  # M86

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 4715.4614  # H89
$ws.Cells.Item(89, 9).Value = 4536.4546  # I89
$ws.Cells.Item(89, 11).Value = 22682.273  # K89
$ws.Cells.Item(89, 13).Value = -17066.273  # M89

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(128, 8).Value = 2392.0  # H128
$ws.Cells.Item(128, 9).Value = 2392.0  # I128
$ws.Cells.Item(128, 11).Value = 7176.0  # K128
$ws.Cells.Item(128, 13).Value = -4686.0  # M128

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3089.8  # H134
$ws.Cells.Item(134, 9).Value = 2700.0  # I134
$ws.Cells.Item(134, 10).Value = 3187.25  # J134
$ws.Cells.Item(134, 11).Value = 8100.0  # K134
$ws.Cells.Item(134, 12).Value = 9561.75  # L134
$ws.Cells.Item(134, 13).Value = -5565.0  # M134
$ws.Cells.Item(134, 14).Value = -14631.75  # N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10424876.0  # H31
$ws.Cells.Item(31, 9).Value = 2643.8572  # I31
$ws.Cells.Item(31, 11).Value = 2643.8572  # K31
$ws.Cells.Item(31, 13).Value = -2348.8572  # M31

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 10424876.0  # H34
$ws.Cells.Item(34, 9).Value = 2643.8572  # I34
$ws.Cells.Item(34, 11).Value = 2643.8572  # K34
$ws.Cells.Item(34, 13).Value = -2441.8572  # M34

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2131.2  # H58
$ws.Cells.Item(58, 10).Value = 3500.0  # J58
$ws.Cells.Item(58, 12).Value = 3500.0  # L58
$ws.Cells.Item(58, 14).Value = -3906.0  # N58

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 3486.25  # H122
$ws.Cells.Item(122, 9).Value = 2933.0  # I122
$ws.Cells.Item(122, 10).Value = 4039.5  # J122
$ws.Cells.Item(122, 11).Value = 8799.0  # K122
$ws.Cells.Item(122, 12).Value = 12118.5  # L122
$ws.Cells.Item(122, 13).Value = -6349.0  # M122
$ws.Cells.Item(122, 14).Value = -17018.5  # N122

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3659.4883  # H132
$ws.Cells.Item(132, 9).Value = 2673.9167  # I132
$ws.Cells.Item(132, 10).Value = 8728.143  # J132
$ws.Cells.Item(132, 11).Value = 8021.750100000001  # K132
$ws.Cells.Item(132, 12).Value = 26184.429  # L132
$ws.Cells.Item(132, 13).Value = -5491.750100000001  # M132
$ws.Cells.Item(132, 14).Value = -31244.429  # N132

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2131.2  # H136
$ws.Cells.Item(136, 10).Value = 3500.0  # J136
$ws.Cells.Item(136, 12).Value = 10500.0  # L136
$ws.Cells.Item(136, 14).Value = -15600.0  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 312.5  # H7
$ws.Cells.Item(7, 9).Value = 350.0  # I7
$ws.Cells.Item(7, 11).Value = 1050.0  # K7
$ws.Cells.Item(7, 13).Value = -938.0  # M7

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 483.9  # H11
$ws.Cells.Item(11, 9).Value = 429.875  # I11
$ws.Cells.Item(11, 11).Value = 1289.625  # K11
$ws.Cells.Item(11, 13).Value = -1149.625  # M11

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 70667.375  # H129
$ws.Cells.Item(129, 10).Value = 94132.336  # J129
$ws.Cells.Item(129, 12).Value = 282397.008  # L129
$ws.Cells.Item(129, 14).Value = -292397.008  # N129

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(136, 8).Value = 919.125  # H136
$ws.Cells.Item(136, 9).Value = 919.125  # I136
$ws.Cells.Item(136, 10).Value = 0.0  # J136
$ws.Cells.Item(136, 11).Value = 2757.375  # K136
$ws.Cells.Item(136, 12).Value = 0.0  # L136
$ws.Cells.Item(136, 13).Value = 2342.625  # M136
$ws.Cells.Item(136, 14).ClearContents()  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(138, 8).Value = 10374.75  # H138
$ws.Cells.Item(138, 9).Value = 7166.3335  # I138
$ws.Cells.Item(138, 11).Value = 21499.0005  # K138
$ws.Cells.Item(138, 13).Value = -16359.0005  # M138

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 0.0  # H22
$ws.Cells.Item(22, 10).Value = 0.0  # J22
$ws.Cells.Item(22, 12).Value = 0.0  # L22
$ws.Cells.Item(22, 14).ClearContents()  # N22

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 90911464.0  # H80
$ws.Cells.Item(80, 9).Value = 250002240.0  # I80
$ws.Cells.Item(80, 11).Value = 250002240.0  # K80
$ws.Cells.Item(80, 13).Value = -250001242.0  # M80

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 90911464.0  # H83
$ws.Cells.Item(83, 9).Value = 250002240.0  # I83
$ws.Cells.Item(83, 11).Value = 1250011200.0  # K83
$ws.Cells.Item(83, 13).Value = -1250006208.0  # M83

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3348015.8  # H122
$ws.Cells.Item(122, 10).Value = 4669.0  # J122
$ws.Cells.Item(122, 12).Value = 14007.0  # L122
$ws.Cells.Item(122, 14).Value = -18907.0  # N122

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 8022.5557  # H126
$ws.Cells.Item(126, 9).Value = 3532.2  # I126
$ws.Cells.Item(126, 11).Value = 10596.6  # K126
$ws.Cells.Item(126, 13).Value = -8126.599999999999  # M126

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2280.8333  # H132
$ws.Cells.Item(132, 9).Value = 2196.3157  # I132
$ws.Cells.Item(132, 10).Value = 2426.818  # J132
$ws.Cells.Item(132, 11).Value = 6588.9471  # K132
$ws.Cells.Item(132, 12).Value = 7280.454000000001  # L132
$ws.Cells.Item(132, 13).Value = -4058.9471  # M132
$ws.Cells.Item(132, 14).Value = -12340.454  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 7998.909  # H132
$ws.Cells.Item(132, 9).Value = 7284.0  # I132
$ws.Cells.Item(132, 11).Value = 21852.0  # K132
$ws.Cells.Item(132, 13).Value = -19322.0  # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 4430.8696  # H136
$ws.Cells.Item(136, 9).Value = 3957.1428  # I136
$ws.Cells.Item(136, 10).Value = 5167.778  # J136
$ws.Cells.Item(136, 11).Value = 11871.4284  # K136
$ws.Cells.Item(136, 12).Value = 15503.334  # L136
$ws.Cells.Item(136, 13).Value = -9321.4284  # M136
$ws.Cells.Item(136, 14).Value = -20603.334  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 11906645.0  # H122
$ws.Cells.Item(122, 10).Value = 31251982.0  # J122
$ws.Cells.Item(122, 12).Value = 93755946.0  # L122
$ws.Cells.Item(122, 14).Value = -93760846.0  # N122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3176.0908  # H132
$ws.Cells.Item(132, 9).Value = 2493.8  # I132
$ws.Cells.Item(132, 11).Value = 7481.400000000001  # K132
$ws.Cells.Item(132, 13).Value = -4951.400000000001  # M132
